# Generate Report for Archive
#
# The localization run moved from "handoff" to "in translation": update the
# Status text everywhere it appears (Overview's per-language status columns
# E/F, and the Status column C on each language sheet), then re-fit those
# columns' widths now that the new text is shorter than the old one.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# 1. Status text: "Ready for handoff" -> "In Translation"
$wsOverview.Range("E2:F3").Value = "In Translation"
$wsZhCn.Range("C2:C3").Value = "In Translation"
$wsDeDe.Range("C2:C3").Value = "In Translation"

# 2. Shrink the now-narrower Status columns to match the shorter text.
$newStatusColWidth = 12.5

$wsOverview.Columns.Item(5).ColumnWidth = $newStatusColWidth
$wsOverview.Columns.Item(6).ColumnWidth = $newStatusColWidth

$wsZhCn.Columns.Item(3).ColumnWidth = $newStatusColWidth

$wsDeDe.Columns.Item(3).ColumnWidth = $newStatusColWidth
